# Apply the TimeReport_FernandoAmado.xlsx update:
#  - correct the Oct 4 ending time (D5)
#  - fill in the two blank rows (9 and 10) with new Development / Debugging entries
#  - move the active cell selection to D11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: fix the recorded ending time ---
$ws.Range("D5").Value = 0.81597222222222221

# --- Row 9: Development entry on 10/24/2025 ---
$ws.Range("A9").Value = 45954

$ws.Range("B9").Value = "Development"

$ws.Range("C9").Value = 0.61111111111111116
$ws.Range("C9").NumberFormat = "h:mm AM/PM"

$ws.Range("D9").Value = 0.89583333333333337
$ws.Range("D9").NumberFormat = "h:mm AM/PM"

$ws.Range("F9").Value = "Worked with group to configure environments for embedded development, build test circuit for experiment with microcontroller, and test adc example with test circuit."

# --- Row 10: Debugging entry on 10/25/2025 ---
$ws.Range("A10").Value = 45955

$ws.Range("B10").Value = "Debugging"

$ws.Range("C10").Value = 0.45833333333333331
$ws.Range("C10").NumberFormat = "h:mm AM/PM"

$ws.Range("D10").Value = 0.66666666666666663
$ws.Range("D10").NumberFormat = "h:mm AM/PM"

$ws.Range("F10").Value = "Debugged personal environment to continue work in github integrated workflow. In the meantime helped set up test circuits and microcontroller experiments."

# --- Recalculate formulas (E column totals, etc.) ---
$excel.Calculate()

# --- Move the active selection to D11 ---
$ws.Range("D11").Select()
